# 10.3.1.xlsx update:
#  - The "section header" row for the Age breakdown (row 19) gets more
#    specific wording ("By age ..." instead of the generic "Age ...").
#  - The "section header" row for the Education breakdown (row 29) gets
#    more specific wording ("By education ..." instead of the generic
#    "Education ...").
# These are the only user-visible content changes; everything else in the
# underlying XML diff (shared-string reindexing) is a mechanical side
# effect of Excel rewriting the shared-strings table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Column A (Kyrgyz) updated first: age section header (row 19), then the
# education section header (row 29).
$ws.Range("A19").Value = "Жаш курагы боюнча (жылдарда)"
$ws.Range("A29").Value = "Билими боюнча"

# Column B (Russian) updated next: age, then education.
$ws.Range("B19").Value = "По возрасту (в годах)"
$ws.Range("B29").Value = "По образованию"

# Column C (English) updated last: age, then education.
$ws.Range("C19").Value = "By age (in years) "
$ws.Range("C29").Value = "By education"
